# Cotações atualizadas - 2025-10-17
# Append a new daily quote row (row 43) below the existing data (rows 2-42).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 43

# Column A: the date, stored as a serial number (2025-10-17 == 45947),
# formatted the same way as the preceding date cells in column A.
$ws.Cells.Item($newRow, 1).Value = 45947
$ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Columns B-E: the fund quotes for the day, stored as text (comma decimal
# separator), matching the format of the existing rows.
$ws.Cells.Item($newRow, 2).Value = "21,7414"
$ws.Cells.Item($newRow, 3).Value = "15,4996"
$ws.Cells.Item($newRow, 4).Value = "15,5859"
$ws.Cells.Item($newRow, 5).Value = "15,5859"
